# Change logic for image slide: images no longer change scale (all squares),
# now laid out in a 2x2-ish grid (two images on top row, one on the
# second row) instead of three images side by side.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$emuPerInch = 914400
$ptPerInch = 72

function EmuToPt($emu) {
    # Convert EMU -> points. A tiny epsilon compensates for floating point
    # truncation in the host's point->EMU re-conversion on write (it floors
    # rather than rounds), so values land back on the exact EMU we want.
    return ($emu / $emuPerInch * $ptPerInch) + 0.00002
}

# Picture 2 (shape index 2)
$pic1 = $s.Shapes.Item(2)
$pic1.Left   = EmuToPt 1891500
$pic1.Top    = EmuToPt 900000
$pic1.Width  = EmuToPt 2493000
$pic1.Height = EmuToPt 2493000

# TextBox 3 - caption for picture 1 (shape index 3)
$cap1 = $s.Shapes.Item(3)
$cap1.Left   = EmuToPt 360000
$cap1.Top    = EmuToPt 3393000
$cap1.Width  = EmuToPt 5556000
$cap1.Height = EmuToPt 216000

# Picture 4 (shape index 4)
$pic2 = $s.Shapes.Item(4)
$pic2.Left   = EmuToPt 7807500
$pic2.Top    = EmuToPt 900000
$pic2.Width  = EmuToPt 2493000
$pic2.Height = EmuToPt 2493000

# TextBox 5 - caption for picture 2 (shape index 5)
$cap2 = $s.Shapes.Item(5)
$cap2.Left   = EmuToPt 6276000
$cap2.Top    = EmuToPt 3393000
$cap2.Width  = EmuToPt 5556000
$cap2.Height = EmuToPt 216000

# Picture 6 (shape index 6)
$pic3 = $s.Shapes.Item(6)
$pic3.Left   = EmuToPt 1891500
$pic3.Top    = EmuToPt 3969000
$pic3.Width  = EmuToPt 2493000
$pic3.Height = EmuToPt 2493000

# TextBox 7 - caption for picture 3 (shape index 7)
$cap3 = $s.Shapes.Item(7)
$cap3.Left   = EmuToPt 360000
$cap3.Top    = EmuToPt 6462000
$cap3.Width  = EmuToPt 5556000
$cap3.Height = EmuToPt 216000
